$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the checklist table with three new "Load button" test cases ---
# Row 30 (C/D/E were blank) now gets the first new test case.
# Rows 31-32 are brand-new rows continuing the numbering/formula pattern.
# Rows 33-34 are brand-new, left blank (same as the former "next" blank row).

# Copy the formatting of row 30 down into the new rows 31:34 first, so the
# new cells inherit the same styles (borders/fill/alignment) as row 30.
$ws.Range("B30:E30").Copy()
$ws.Range("B31:E31").PasteSpecial(-4122)
$ws.Range("B32:E32").PasteSpecial(-4122)
$ws.Range("B33:E34").PasteSpecial(-4122)

# Continue the running Tc_N numbering formula into the two populated rows.
$ws.Range("B31").Formula = "=B30+1"
$ws.Range("B32").Formula = "=B31+1"

# Row 30: Load button / Verify load button appear in the ui
$ws.Range("C30").Value = "Load button"
$ws.Range("D30").Value = "Verify load button appear in the ui"
$ws.Range("E30").Value = "Ready"

# Row 31: Load button / verify when click on load button the number of movies in the list increase
$ws.Range("C31").Value = "Load button"
$ws.Range("D31").Value = "verify when click on load button the number of movies in the list increase"
$ws.Range("E31").Value = "Ready"

# Row 32: Load button / Verify click in the load button 3 times
$ws.Range("C32").Value = "Load button"
$ws.Range("D32").Value = "Verify click in the load button 3 times "
$ws.Range("E32").Value = "Ready"

# Rows 33 and 34 stay blank (format only, already applied above).

# --- Extend the validation dropdown to cover the new rows ---
$ws.Range("E3:E34").Validation.Delete()
$ws.Range("E3:E34").Validation.Add(3, 1, 1, '"Pass,Fail,Ready"')

# --- Update the view: scroll position, zoom, and selection ---
$ws.Range("B33:B34").Select()
$excel.ActiveWindow.Zoom = 76
